# "Write script for use excel Data."
# MyData.xlsx / Sheet1: two of the "Name" column values were retyped.
# Row 4's Name ("C") becomes "Auntor" and row 5's Name ("D") becomes "Acharja".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = "Auntor"
$ws.Range("B5").Value = "Acharja"

# Leave the cell cursor on C11, matching where the editor's selection ended up.
$ws.Range("C11").Select()
